# Generate Report for Handback
# Refreshes the handback-status report with newly computed handoff/handback
# timestamps (and the corresponding correspond-file values) for the
# 63679998-18cd-485d-97e6-76d9b68749bd file, across the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-26 10:46:31"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("G3").Value = "2016-08-26 10:46:27"
$zhcn.Range("H3").Value = "2016-08-26 10:46:44"
$zhcn.Range("J3").Value = "2016-08-26 10:46:27"
$zhcn.Range("K3").Value = "63679998-18cd-485d-97e6-76d9b68749bd.134fa91a48401a5a341eeb3756855d81b3abc47e.zh-cn.xlf"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-26 10:46:31"
$dede.Range("K3").Value = "2016-08-26 10:46:50"
